$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.955.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.03%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.116.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.95%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''577.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.47%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''172.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +2.32%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.06%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -0.65%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''6.44'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -2.49%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -1.08%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.484'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.37%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -1.39%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''37.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.48%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  -1.31%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.634.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +1.06%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''66.889.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.02%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -0.29%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.118.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +1.03%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''16.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +0.74%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''476.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +2.09%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = '''Polygon'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '''https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '''0.711'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.56%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = '''Uniswap'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '''7.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +5.63%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = '''Litecoin'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = '''83.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.74%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = '''InternetComputer(DFINITY)'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = '''13.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +2.87%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -3.47%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''10.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.64%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -0.06%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''7.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.42%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -1.36%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +0.45%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''28.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +1.13%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +0.37%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.0₃0944'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -7.82%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +0.03%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.76%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.975'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -2.94%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''46.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.35%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''50.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.09%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''2.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -2.76%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.314'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -2.06%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.123'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.36%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''8.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -0.21%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''2.817.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.77%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''383.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -0.12%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -1.88%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -10.26%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''135.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.20%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = '''24.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +0.00%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -1.70%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -0.75%  '
$ws.Range("E51").Style = "Normal"
